$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,6,10,12,13,14,15,18,19,20,21,22,24,26,28,29,32,36,38,39,40,41,44,45,46,47,48,50,52,54,55,58,62,64,65,66,67,70,71,72,73,74,76,78,83,84,85,86,90,92,99,101,109,110,111,112,116,118,125,127,135,136,137,138,142,144,151,153)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $val = $cell.Value2
    if ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value2 = "system, backup@backdoor.com, System"
    } elseif ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
}
